$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.397.11"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.463.42"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "576.24"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "160.72"
$ws.Range("E6").Value = "  +3.94%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.465.54"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.85%  "
$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").Value = "4.054.33"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("E15").Value = "  +5.44%  "
$ws.Range("D16").Value = "29.13"
$ws.Range("E16").Value = "  +7.07%  "
$ws.Range("D17").Value = "64.441.93"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "3.460.92"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").Value = "6.45"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "14.53"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").Value = "387.49"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "8.27"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "73.32"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +10.41%  "
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  +9.21%  "
$ws.Range("D32").Value = "2.03"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "6.65"
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "7.12"
$ws.Range("E36").Value = "  +5.23%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "160.94"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "27.59"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "2.925.20"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "0.0322"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "4.53"
$ws.Range("E44").Value = "  +4.96%  "
$ws.Range("D45").Value = "42.56"
$ws.Range("E45").Value = "  +3.48%  "
$ws.Range("D46").Value = "0.773"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "23.99"
$ws.Range("E47").Value = "  +8.38%  "
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("D49").Value = "2.21"
$ws.Range("E49").Value = "  +15.45%  "
$ws.Range("D50").Value = "0.109"
$ws.Range("E50").Value = "  +4.81%  "
$ws.Range("D51").Value = "6.61"
$ws.Range("E51").Value = "  +4.55%  "
